$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old A1:E6 block entirely (rows/cols are moving to B8:F14)
$ws.Range("A1:E6").Clear()

# Data rows, shifted from A1:E6 to B9:F14
$data = @(
    @(1000, 1, "Roolers1", "Some description1", "NoImage.jpg"),
    @(1001, 5, "Roolers2", "Some description2", "NoImage.jpg"),
    @(1002, 2, "Roolers3", "Some description3", "NoImage.jpg"),
    @("hgfjf", 4, "Roolers4", "Some description4", "NoImage.jpg"),
    @(1004, 3, "Roolers5", "Some description5", "NoImage.jpg"),
    @(1005, 5, "Roolers6", "Some description6", "NoImage.jpg")
)

$row = 9
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $ws.Cells.Item($row, 5).Value = $item[3]
    $ws.Cells.Item($row, 6).Value = $item[4]
    $row++
}

# New row 8: summary / new item being added, with gaps at D8 and E8
$ws.Cells.Item(8, 2).Value = 999
$ws.Cells.Item(8, 6).Value = "sfvdf"
$ws.Cells.Item(8, 3).Value = "фитнесс"

$ws.Range("E9").Select()
